# Update timestamps on the "Generate Report for Handback" regeneration.
$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date for the 7188329f... entry (row 3)
$overview.Range("G3").Value = "2016-08-30 22:49:20"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for the 7188329f... entry (row 3)
$zhcn.Range("H3").Value = "2016-08-30 22:49:15"
$zhcn.Range("K3").Value = "2016-08-30 22:49:34"

# de-de sheet: Correspond Handback DateTime for the 7188329f... entry (row 3)
$dede.Range("K3").Value = "2016-08-30 22:49:42"
